# fixed some bugs in genSymbolVals2
# Reorders the data rows (rows 2-23, columns A-F) on Sheet1 to their
# corrected values/positions per the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(101, 9, 30, 15, 60, 15)
    3  = @(1001, 18, 30, 75, 60, 72)
    4  = @(1202, 2, 10, 10, 10, 10)
    5  = @(701, 3, 90, 45, 97, 15)
    6  = @(201, 9, 30, 15, 45, 30)
    7  = @(801, 3, 67, 65, 52, 45)
    8  = @(1203, 3, 15, 15, 15, 15)
    9  = @(901, 16, 15, 45, 60, 60)
    10 = @(301, 6, 45, 30, 60, 45)
    11 = @(501, 9, 52, 30, 75, 45)
    12 = @(601, 9, 60, 67, 60, 42)
    13 = @(902, 1, 0, 0, 0, 0)
    14 = @(401, 9, 48, 67, 75, 45)
    15 = @(1201, 2, 10, 10, 10, 10)
    16 = @(1101, 0, 15, 30, 30, 0)
    17 = @(802, 0, 4, 5, 4, 0)
    18 = @(1, 0, 2, 2, 2, 2)
    19 = @(3, 0, 3, 3, 3, 3)
    20 = @(2, 0, 2, 2, 2, 2)
    21 = @(502, 0, 4, 0, 0, 0)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
